# Updates the cryptos list (coin name / link / price / 1h volume %) per the
# scraped coinranking.com data refresh for this run of the GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map: row number -> hashtable of column letter -> new value (as text)
$changes = @{
    2  = @{ D = "66.208.45";   E = "  -2.25%  " }
    3  = @{ D = "3.834.09";    E = "  +1.94%  " }
    4  = @{                    E = "  -0.36%  " }
    5  = @{ D = "422.05";      E = "  +0.16%  " }
    6  = @{ D = "127.54";      E = "  -3.53%  " }
    7  = @{ D = "3.831.45";    E = "  +2.28%  " }
    8  = @{ D = "0.602";       E = "  -7.56%  " }
    10 = @{ D = "0.715";       E = "  -7.68%  " }
    11 = @{                    E = "  -12.55%  " }
    12 = @{ D = "0.0000345";   E = "  -19.69%  " }
    13 = @{ D = "39.99";       E = "  -6.71%  " }
    14 = @{ D = "4.418.11";    E = "  +1.24%  " }
    15 = @{ D = "9.89";        E = "  -5.07%  " }
    16 = @{ D = "15.80";       E = "  +21.03%  " }
    17 = @{ B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "3.832.61"; E = "  +2.03%  " }
    18 = @{ B = "TRON";         C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx";         D = "0.137";    E = "  -1.85%  " }
    19 = @{ D = "19.41";       E = "  -5.87%  " }
    20 = @{ D = "66.376.21";   E = "  -1.99%  " }
    21 = @{                    E = "  -6.85%  " }
    22 = @{ D = "401.38";      E = "  -11.05%  " }
    23 = @{ D = "14.21";       E = "  -11.07%  " }
    24 = @{ D = "83.48";       E = "  -7.31%  " }
    25 = @{ D = "2.97";        E = "  -3.99%  " }
    26 = @{ B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "36.78"; E = "  -4.29%  " }
    27 = @{ B = "LEO";             C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";             D = "5.76";  E = "  +12.84%  " }
    28 = @{                    E = "  -5.17%  " }
    29 = @{ D = "9.34";        E = "  -7.80%  " }
    30 = @{ D = "699.60";      E = "  +1.77%  " }
    31 = @{ B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "2.74";   E = "  -0.68%  " }
    32 = @{ B = "Hedera";  C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.120";   E = "  -4.36%  " }
    33 = @{ D = "12.21";       E = "  -4.12%  " }
    34 = @{ D = "7.46";        E = "  +3.41%  " }
    35 = @{                    E = "  -10.53%  " }
    36 = @{ B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D = "37.61"; E = "  -10.52%  " }
    37 = @{ B = "Dai";               C = "https://coinranking.com/coin/MoTuySvg7+dai-dai";               D = "1.00";  E = "  +0.15%  " }
    38 = @{ D = "54.76";       E = "  -4.37%  " }
    39 = @{ D = "0.0₃0756";     E = "  +0.49%  " }
    40 = @{ D = "0.0448";      E = "  -9.36%  " }
    41 = @{                    E = "  -5.00%  " }
    42 = @{                    E = "  +0.26%  " }
    43 = @{                    E = "  -9.98%  " }
    44 = @{ D = "4.42";        E = "  +1.60%  " }
    45 = @{ D = "3.28";        E = "  -3.47%  " }
    46 = @{ D = "143.49";      E = "  -2.76%  " }
    47 = @{ D = "3.08";        E = "  -2.29%  " }
    48 = @{ B = "ARBITRUM";   C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D = "2.05"; E = "  -4.13%  " }
    49 = @{ B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "25.86"; E = "  -7.36%  " }
    50 = @{                    E = "  -4.88%  " }
    51 = @{ D = "2.71";        E = "  -7.24%  " }
}

# These "Price" cells look like plain numbers/decimals (e.g. "422.05").
# The source data stores them as plain text, so when writing them back we
# prefix with a leading apostrophe to stop Excel from reinterpreting them
# as numeric values (which would also reformat them, e.g. drop trailing
# zeros or add float noise).
function Test-NumericText([string]$s) {
    return $s -match '^[+-]?[0-9]+(\.[0-9]+)?$'
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $value = $cols[$col]
        if ($col -eq "D" -and (Test-NumericText $value)) {
            $ws.Range($addr).Value = "'" + $value
        } else {
            $ws.Range($addr).Value = $value
        }
    }
}
